# Updates the "cryptos" sheet with refreshed prices / 1h volume changes,
# and fixes the ordering of a few coin rows (PancakeSwap/PEPE, OKB/Bittensor,
# WhiteBITCoin/Mantle) whose Coin/Link/Price/Volume data had been swapped.
# Numeric-looking price strings are forced to Text format before assignment
# so Excel keeps them as text (matching the source data) instead of
# auto-converting them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.759.86'
$ws.Range("E2").Value = '  -4.66%  '

$ws.Range("D3").Value = '2.477.32'
$ws.Range("E3").Value = '  -3.52%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '531.83'
$ws.Range("E5").Value = '  -2.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.33'
$ws.Range("E6").Value = '  -7.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.34%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  -4.99%  '

$ws.Range("D9").Value = '2.501.39'
$ws.Range("E9").Value = '  -2.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0992'
$ws.Range("E10").Value = '  -4.61%  '

$ws.Range("E11").Value = '  -2.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.46'
$ws.Range("E12").Value = '  +0.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.348'
$ws.Range("E13").Value = '  -4.33%  '

$ws.Range("D14").Value = '2.905.36'
$ws.Range("E14").Value = '  -3.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.62'
$ws.Range("E15").Value = '  -7.03%  '

$ws.Range("D16").Value = '58.624.94'
$ws.Range("E16").Value = '  -4.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("E17").Value = '  -4.49%  '

$ws.Range("D18").Value = '2.478.40'
$ws.Range("E18").Value = '  -3.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.25'
$ws.Range("E19").Value = '  -2.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.24'
$ws.Range("E20").Value = '  -6.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.01'
$ws.Range("E21").Value = '  -4.82%  '

$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.70'
$ws.Range("E23").Value = '  -5.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.48'
$ws.Range("E24").Value = '  -4.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.431'
$ws.Range("E25").Value = '  -12.67%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.160'
$ws.Range("E27").Value = '  -4.99%  '

$ws.Range("D28").Value = '2.576.74'
$ws.Range("E28").Value = '  -4.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.69'
$ws.Range("E29").Value = '  -5.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.74'
$ws.Range("E30").Value = '  -8.87%  '

$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0₃0756'
$ws.Range("E31").Value = '  -9.39%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.77'
$ws.Range("E32").Value = '  -6.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.23'
$ws.Range("E33").Value = '  -4.34%  '

$ws.Range("E34").Value = '  -0.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.75'
$ws.Range("E35").Value = '  -2.05%  '

$ws.Range("E36").Value = '  -0.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.43'
$ws.Range("E37").Value = '  -3.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.31'
$ws.Range("E38").Value = '  -8.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.59'
$ws.Range("E39").Value = '  -12.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.78'
$ws.Range("E40").Value = '  -2.65%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '304.00'
$ws.Range("E41").Value = '  -8.57%  '

$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.60'
$ws.Range("E42").Value = '  -2.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.62'
$ws.Range("E43").Value = '  -8.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.778'
$ws.Range("E44").Value = '  -16.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.994'
$ws.Range("E45").Value = '  -0.57%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.592'
$ws.Range("E46").Value = '  -1.86%  '

$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.76'
$ws.Range("E47").Value = '  -1.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.37'
$ws.Range("E48").Value = '  +1.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0921'
$ws.Range("E49").Value = '  -4.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.41'
$ws.Range("E50").Value = '  -5.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0512'
$ws.Range("E51").Value = '  -6.09%  '
